$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Fill in the "Model" column (B) for BNC sockets and Barrel connector rows,
# pointing readers to the updated detailed BOM CSV (fixes broken component links).
$ws.Range("B7").Value = "see detailed_bom.csv"
$ws.Range("B8").Value = "see detailed_bom.csv"

# Update the active selection to match the saved view state.
$ws.Range("B8").Select()
